# [Feat 2269] Add support of parameters worksheet metadata.
# Rename the "PARAMETERS TODO" sheet to "PARAMETERS" and populate its header
# row with the ACTION / TC_OWNER_PATH / TC_OWNER_ID / TC_PARAM_ID /
# TC_PARAM_NAME / TC_PARAM_DESCRIPTION columns, then make it the active sheet.

$wb = $excel.ActiveWorkbook

# Rename "PARAMETERS TODO" -> "PARAMETERS"
$ws = $wb.Worksheets.Item("PARAMETERS TODO")
$ws.Name = "PARAMETERS"

# Header row for the PARAMETERS worksheet.
$ws.Range("A1").Value = "ACTION"
$ws.Range("B1").Value = "TC_OWNER_PATH"
$ws.Range("C1").Value = "TC_OWNER_ID"
$ws.Range("D1").Value = "TC_PARAM_ID"
$ws.Range("E1").Value = "TC_PARAM_NAME"
$ws.Range("F1").Value = "TC_PARAM_DESCRIPTION"

# Size the new columns to fit their header text.
$ws.Columns.Item(2).ColumnWidth = 16.0
$ws.Columns.Item(3).ColumnWidth = 13.0
$ws.Columns.Item(4).ColumnWidth = 12.666666
$ws.Columns.Item(5).ColumnWidth = 16.5
$ws.Columns.Item(6).ColumnWidth = 22.833333

# Make PARAMETERS the active sheet/tab and leave the selection on B19.
$ws.Activate()
$ws.Range("B19").Select()
